$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (was 45406 / 24-Apr-2024, now 45436 / 24-May-2024)
$ws.Range("A1").Value = Get-Date -Year 2024 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Update price values
$ws.Range("D22").Value = 416.227
$ws.Range("D23").Value = 527.458
$ws.Range("D34").Value = 410.487
$ws.Range("D35").Value = 567.287
$ws.Range("D45").Value = 503.776
$ws.Range("D46").Value = 565.133

